# Readying the blog post
# - Add a new "within_universe" worksheet with character/dressed_as/season/episode data
# - Fix a couple of swapped costume_detail / costume_category values for Dwight (rows 94-95)
# - Restore the view state (active cell selections) on both sheets

$wb = $excel.ActiveWorkbook
$rawSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Swap the mixed-up costume_detail (D) / costume_category (F) values for
#    the two Dwight "Pig" / "Jack-o-lantern" rows.
# ---------------------------------------------------------------------------
$d94 = $rawSheet.Range("D94").Value2
$f94 = $rawSheet.Range("F94").Value2
$d95 = $rawSheet.Range("D95").Value2
$f95 = $rawSheet.Range("F95").Value2

$rawSheet.Range("D94").Value2 = $d95
$rawSheet.Range("F94").Value2 = $f95
$rawSheet.Range("D95").Value2 = $d94
$rawSheet.Range("F95").Value2 = $f94

# ---------------------------------------------------------------------------
# 2. Add the new "within_universe" sheet right after office_costumes_raw.
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $rawSheet)
$newSheet.Name = "within_universe"

$data = @(
    @("character", "dressed_as", "season", "episode"),
    @("Creed",   "Creed",    9, 5),
    @("Dwight",  "Jim",      3, 21),
    @("Dwight",  "Toby",     7, 6),
    @("Dwight",  "Kevin",    7, 11),
    @("Dwight",  "Meredith", 7, 11),
    @("Dwight",  "Pam",      7, 11),
    @("Dwight",  "Stanley",  7, 11),
    @("Jim",     "Dwight",   3, 21),
    @("Michael", "Darryl",   7, 6),
    @("Michael", "Angela",   7, 21),
    @("Michael", "Jim",      7, 21),
    @("Michael", "Jo",       7, 21),
    @("Michael", "Phyllis",  7, 21),
    @("Nellie",  "Toby",     9, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 1
    $rowData = $data[$i]
    $newSheet.Cells.Item($rowNum, 1).Value2 = $rowData[0]
    $newSheet.Cells.Item($rowNum, 2).Value2 = $rowData[1]
    $newSheet.Cells.Item($rowNum, 3).Value2 = $rowData[2]
    $newSheet.Cells.Item($rowNum, 4).Value2 = $rowData[3]
}

# Leave the cursor below the data on the new sheet, as in the authored file.
[void]$newSheet.Range("B16").Select()

# ---------------------------------------------------------------------------
# 3. Restore focus / scroll state on the raw sheet.
# ---------------------------------------------------------------------------
[void]$rawSheet.Activate()
[void]$rawSheet.Range("D96").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 79
$win.ScrollColumn = 1
